$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update elevation type text in B2
$ws.Range("B2").Value = "n"

# Fill in the previously blank computed values on the "Total Cost ($)" row (row 3)
$ws.Range("O3").Value = 1024.8
$ws.Range("P3").Value = 559.8
$ws.Range("Q3").Value = 194.4
